$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.530.86"
$ws.Range("E2").Value = "  +4.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.000.06"
$ws.Range("E3").Value = "  +4.85%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.19"
$ws.Range("E5").Value = "  +9.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.97"
$ws.Range("E6").Value = "  +11.07%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +8.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.56"
$ws.Range("E9").Value = "  +15.43%  "

$ws.Range("E10").Value = "  +14.60%  "

$ws.Range("E12").Value = "  +6.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.515.37"
$ws.Range("E13").Value = "  +4.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.44"
$ws.Range("E14").Value = "  +11.65%  "

$ws.Range("E15").Value = "  +17.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.562.21"
$ws.Range("E16").Value = "  +4.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.005.79"
$ws.Range("E17").Value = "  +5.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.78"
$ws.Range("E18").Value = "  +10.31%  "

$ws.Range("E19").Value = "  +11.10%  "

$ws.Range("E20").Value = "  +12.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.45"
$ws.Range("E21").Value = "  +12.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("E23").Value = "  +10.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.39"
$ws.Range("E24").Value = "  +8.05%  "

$ws.Range("E25").Value = "  +14.56%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  +16.33%  "

$ws.Range("E28").Value = "  +9.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.04"
$ws.Range("E29").Value = "  +16.78%  "

$ws.Range("E30").Value = "  +16.48%  "

$ws.Range("E31").Value = "  +12.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.61"
$ws.Range("E32").Value = "  +12.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.64"
$ws.Range("E33").Value = "  +13.41%  "

$ws.Range("E34").Value = "  +10.42%  "

$ws.Range("E35").Value = "  +6.10%  "

$ws.Range("E36").Value = "  +6.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0675"
$ws.Range("E37").Value = "  +11.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.91"
$ws.Range("E38").Value = "  +5.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.034.03"
$ws.Range("E39").Value = "  +5.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.65"
$ws.Range("E40").Value = "  +4.55%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.649"
$ws.Range("E42").Value = "  +8.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.259.53"
$ws.Range("E43").Value = "  +12.02%  "

$ws.Range("E44").Value = "  +8.26%  "

$ws.Range("E45").Value = "  +9.43%  "

$ws.Range("E46").Value = "  +8.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.98"
$ws.Range("E47").Value = "  +26.88%  "

$ws.Range("E48").Value = "  +12.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.74"
$ws.Range("E49").Value = "  +9.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.98"
$ws.Range("E50").Value = "  +8.90%  "

$ws.Range("E51").Value = "  +12.24%  "
